# Generate Report for Handoff
# Updates the localization-status report: the translation run that was
# "In Translation" has now completed and is ready to hand off to the next
# stage, and the handoff timestamps advance a few seconds to when the
# report was (re)generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 00:38:53"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 00:38:49"

# --- de-de sheet ------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 00:38:53"

# --- widen the Status / Latest Handoff Datetime columns to fit the new,
#     longer "Ready for handoff" text (mirrors the column autosize that
#     happens when the report is regenerated) -------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
